$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force text storage for the Price/Volume columns so purely-numeric-looking
# strings (e.g. "150.98") are not auto-coerced into floating point numbers.
$priceVolRange = $ws.Range("D2:E51")
$priceVolRange.NumberFormat = "@"

$ws.Range('D2').Value = '62.394.73'
$ws.Range('E2').Value = '  +1.90%  '
$ws.Range('D3').Value = '3.001.53'
$ws.Range('E3').Value = '  +0.85%  '
$ws.Range('E4').Value = '  +0.11%  '
$ws.Range('D5').Value = '592.89'
$ws.Range('E5').Value = '  -0.36%  '
$ws.Range('D6').Value = '150.98'
$ws.Range('E6').Value = '  +5.51%  '
$ws.Range('E7').Value = '  +0.08%  '
$ws.Range('D8').Value = '2.993.31'
$ws.Range('E8').Value = '  +0.62%  '
$ws.Range('D9').Value = '0.508'
$ws.Range('E9').Value = '  -0.69%  '
$ws.Range('D10').Value = '6.81'
$ws.Range('E10').Value = '  +12.96%  '
$ws.Range('D11').Value = '0.148'
$ws.Range('E11').Value = '  +0.00%  '
$ws.Range('D12').Value = '0.459'
$ws.Range('E12').Value = '  +1.94%  '
$ws.Range('D13').Value = '0.0000230'
$ws.Range('E13').Value = '  +1.81%  '
$ws.Range('D14').Value = '35.20'
$ws.Range('E14').Value = '  +3.55%  '
$ws.Range('E15').Value = '  -0.20%  '
$ws.Range('D16').Value = '3.502.54'
$ws.Range('E16').Value = '  +1.01%  '
$ws.Range('B17').Value = 'Polkadot'
$ws.Range('C17').Value = 'https://coinranking.com/coin/25W7FG7om+polkadot-dot'
$ws.Range('D17').Value = '7.01'
$ws.Range('E17').Value = '  +1.78%  '
$ws.Range('B18').Value = 'WrappedBTC'
$ws.Range('C18').Value = 'https://coinranking.com/coin/x4WXHge-vvFY+wrappedbtc-wbtc'
$ws.Range('D18').Value = '62.407.53'
$ws.Range('E18').Value = '  +1.99%  '
$ws.Range('D19').Value = '3.008.21'
$ws.Range('E19').Value = '  +1.06%  '
$ws.Range('D20').Value = '444.49'
$ws.Range('E20').Value = '  -0.06%  '
$ws.Range('D21').Value = '14.10'
$ws.Range('E21').Value = '  +1.56%  '
$ws.Range('D22').Value = '0.690'
$ws.Range('E22').Value = '  +1.62%  '
$ws.Range('D23').Value = '7.45'
$ws.Range('E23').Value = '  +1.82%  '
$ws.Range('D24').Value = '82.27'
$ws.Range('E24').Value = '  +1.59%  '
$ws.Range('D25').Value = '11.12'
$ws.Range('E25').Value = '  +3.90%  '
$ws.Range('D26').Value = '2.25'
$ws.Range('E26').Value = '  +3.59%  '
$ws.Range('D27').Value = '12.22'
$ws.Range('E27').Value = '  +2.07%  '
$ws.Range('E28').Value = '  -0.01%  '
$ws.Range('D29').Value = '7.43'
$ws.Range('E29').Value = '  +3.84%  '
$ws.Range('D30').Value = '2.25'
$ws.Range('E30').Value = '  +10.28%  '
$ws.Range('D31').Value = '2.67'
$ws.Range('E31').Value = '  -0.68%  '
$ws.Range('E32').Value = '  +0.11%  '
$ws.Range('D33').Value = '27.37'
$ws.Range('E33').Value = '  +0.99%  '
$ws.Range('D34').Value = '0.109'
$ws.Range('E34').Value = '  +0.45%  '
$ws.Range('D35').Value = '0.0₃0857'
$ws.Range('E35').Value = '  +6.93%  '
$ws.Range('D36').Value = '1.03'
$ws.Range('E36').Value = '  +1.85%  '
$ws.Range('D37').Value = '5.83'
$ws.Range('E37').Value = '  +1.55%  '
$ws.Range('D38').Value = '3.08'
$ws.Range('E38').Value = '  +9.64%  '
$ws.Range('D39').Value = '2.08'
$ws.Range('E39').Value = '  +3.95%  '
$ws.Range('B40').Value = 'OKB'
$ws.Range('C40').Value = 'https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb'
$ws.Range('D40').Value = '50.31'
$ws.Range('E40').Value = '  +0.44%  '
$ws.Range('B41').Value = 'Kaspa'
$ws.Range('C41').Value = 'https://coinranking.com/coin/V8GxkwWow+kaspa-kas'
$ws.Range('D41').Value = '0.127'
$ws.Range('E41').Value = '  +3.69%  '
$ws.Range('D42').Value = '8.94'
$ws.Range('E42').Value = '  +0.24%  '
$ws.Range('D43').Value = '44.12'
$ws.Range('E43').Value = '  +12.22%  '
$ws.Range('D44').Value = '0.302'
$ws.Range('E44').Value = '  +12.57%  '
$ws.Range('D45').Value = '0.0356'
$ws.Range('E45').Value = '  +2.55%  '
$ws.Range('D46').Value = '383.44'
$ws.Range('E46').Value = '  -0.28%  '
$ws.Range('D47').Value = '2.687.97'
$ws.Range('E47').Value = '  +0.41%  '
$ws.Range('D48').Value = '133.23'
$ws.Range('E48').Value = '  +2.26%  '
$ws.Range('D49').Value = '26.16'
$ws.Range('E49').Value = '  +12.57%  '
$ws.Range('D51').Value = '2.25'
$ws.Range('E51').Value = '  +5.32%  '

# Clear the temporary text format so cells keep their original (default) style.
$priceVolRange.ClearFormats()
